$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D stay text (matches original inlineStr cells)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.819.43"
$ws.Range("E2").Value = "  -1.29%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.324.02"
$ws.Range("E3").Value = "  +0.79%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.75"
$ws.Range("E5").Value = "  -1.95%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.43"
$ws.Range("E6").Value = "  -2.89%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.507"
$ws.Range("E7").Value = "  -4.59%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.504"
$ws.Range("E9").Value = "  -4.32%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.45"
$ws.Range("E10").Value = "  -5.79%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.14"
$ws.Range("E11").Value = "  -0.22%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0793"
$ws.Range("E12").Value = "  -2.16%  "

$ws.Range("E13").Value = "  +0.64%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.75"
$ws.Range("E14").Value = "  -3.96%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.67"
$ws.Range("E15").Value = "  +4.00%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.349.54"
$ws.Range("E16").Value = "  +1.76%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.818"
$ws.Range("E17").Value = "  +1.17%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.734.75"
$ws.Range("E18").Value = "  -1.22%  "

$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0903"
$ws.Range("E19").Value = "  -2.39%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.13"
$ws.Range("E20").Value = "  -0.53%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.58"
$ws.Range("E21").Value = "  -4.98%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.19"
$ws.Range("E22").Value = "  +1.58%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.51"
$ws.Range("E23").Value = "  -2.91%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.99"
$ws.Range("E24").Value = "  -1.78%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.53"
$ws.Range("E25").Value = "  -3.02%  "

$ws.Range("E26").Value = "  +0.13%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.33"
$ws.Range("E27").Value = "  +1.89%  "

$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.97"
$ws.Range("E28").Value = "  -0.40%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.26"
$ws.Range("E29").Value = "  -1.80%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.75"
$ws.Range("E30").Value = "  -5.89%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.21"
$ws.Range("E31").Value = "  -4.42%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "160.59"
$ws.Range("E32").Value = "  -4.19%  "

$ws.Range("E33").Value = "  +0.03%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.05"
$ws.Range("E34").Value = "  -4.26%  "

$ws.Range("E35").Value = "  -3.27%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.58"
$ws.Range("E36").Value = "  +2.79%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0718"
$ws.Range("E37").Value = "  -3.30%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "16.99"
$ws.Range("E38").Value = "  -6.54%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.88"
$ws.Range("E39").Value = "  -5.51%  "

$ws.Range("E40").Value = "  -2.76%  "

$ws.Range("E41").Value = "  -4.79%  "

$ws.Range("E42").Value = "  -3.33%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.55"
$ws.Range("E43").Value = "  -5.33%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.002.92"
$ws.Range("E44").Value = "  +1.22%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0282"
$ws.Range("E45").Value = "  -4.07%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "18.76"
$ws.Range("E46").Value = "  -1.72%  "

$ws.Range("E47").Value = "  +1.48%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.88"
$ws.Range("E48").Value = "  -4.10%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "55.27"
$ws.Range("E49").Value = "  -1.10%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.87"
$ws.Range("E50").Value = "  -2.04%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.551.74"
$ws.Range("E51").Value = "  +0.78%  "
